# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E2) and "Correspond Handback DateTime" (H2)
# timestamps on the zh-cn and de-de worksheets to reflect the new report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 02:50:59"
$wsZhCn.Range("H2").Value = "2016-03-20 02:51:18"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 02:51:03"
$wsDeDe.Range("H2").Value = "2016-03-20 02:51:24"
